$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2226.6843
$ws.Range("I19").Value = 3733.4707
$ws.Range("J19").Value = 1006.9048
$ws.Range("K19").Value = 3733.4707
$ws.Range("L19").Value = 1006.9048
$ws.Range("M19").Value = -3558.4707
$ws.Range("N19").Value = -1356.9048
$ws.Range("H20").Value = 7705.25
$ws.Range("I20").Value = 940.3333
$ws.Range("J20").Value = 28000
$ws.Range("K20").Value = 940.3333
$ws.Range("L20").Value = 28000
$ws.Range("M20").Value = -710.3333
$ws.Range("N20").Value = -28460
$ws.Range("H34").Value = 27999.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 27999.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 27999.5
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -28405.5
$ws.Range("H35").Value = 7705.25
$ws.Range("I35").Value = 940.3333
$ws.Range("J35").Value = 28000
$ws.Range("K35").Value = 940.3333
$ws.Range("L35").Value = 28000
$ws.Range("M35").Value = -561.3333
$ws.Range("N35").Value = -28758
$ws.Range("H36").Value = 27999.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 27999.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 27999.5
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -29429.5
$ws.Range("H44").Value = 28000
$ws.Range("J44").Value = 28000
$ws.Range("L44").Value = 28000
$ws.Range("N44").Value = -28924
$ws.Range("H47").Value = 25000
$ws.Range("J47").Value = 25000
$ws.Range("L47").Value = 25000
$ws.Range("N47").Value = -26944
$ws.Range("H51").Value = 5194.607
$ws.Range("I51").Value = 3542.9285
$ws.Range("J51").Value = 6846.2856
$ws.Range("K51").Value = 3542.9285
$ws.Range("L51").Value = 6846.2856
$ws.Range("M51").Value = -3058.9285
$ws.Range("N51").Value = -7814.2856
$ws.Range("H70").Value = 2680.5483
$ws.Range("I70").Value = 4369.3335
$ws.Range("J70").Value = 1097.3125
$ws.Range("K70").Value = 13108.0005
$ws.Range("L70").Value = 3291.9375
$ws.Range("M70").Value = -12838.0005
$ws.Range("N70").Value = -3831.9375
$ws.Range("H73").Value = 2680.5483
$ws.Range("I73").Value = 4369.3335
$ws.Range("J73").Value = 1097.3125
$ws.Range("K73").Value = 13108.0005
$ws.Range("L73").Value = 3291.9375
$ws.Range("M73").Value = -12172.0005
$ws.Range("N73").Value = -5163.9375
$ws.Range("H100").Value = 4497.0967
$ws.Range("I100").Value = 1624.3636
$ws.Range("J100").Value = 6077.1
$ws.Range("K100").Value = 1624.3636
$ws.Range("L100").Value = 6077.1
$ws.Range("M100").Value = -1083.3636
$ws.Range("N100").Value = -7159.1
$ws.Range("H132").Value = 4837.2144
$ws.Range("I132").Value = 1165.6666
$ws.Range("J132").Value = 11446
$ws.Range("K132").Value = 3496.9998
$ws.Range("L132").Value = 34338
$ws.Range("M132").Value = -966.9998000000001
$ws.Range("N132").Value = -39398
$ws.Range("H137").Value = 4424.64
$ws.Range("I137").Value = 5051.143
$ws.Range("J137").Value = 3627.2727
$ws.Range("K137").Value = 15153.429
$ws.Range("L137").Value = 10881.8181
$ws.Range("M137").Value = -12603.429
$ws.Range("N137").Value = -15981.8181
$ws.Range("H141").Value = 1030
$ws.Range("I141").Value = 1030
$ws.Range("K141").Value = 3090
$ws.Range("M141").Value = 2090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5377
$ws.Range("I102").Value = 1700.75
$ws.Range("J102").Value = 14200
$ws.Range("K102").Value = 1700.75
$ws.Range("L102").Value = 14200
$ws.Range("M102").Value = -78.75
$ws.Range("N102").Value = -17444
$ws.Range("H132").Value = 28002.324
$ws.Range("I132").Value = 35215.902
$ws.Range("K132").Value = 105647.706
$ws.Range("M132").Value = -103117.706

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 49906
$ws.Range("J95").Value = 49906
$ws.Range("L95").Value = 49906
$ws.Range("N95").Value = -55398
$ws.Range("H97").Value = 6864.778
$ws.Range("I97").Value = 4097.875
$ws.Range("J97").Value = 29000
$ws.Range("K97").Value = 4097.875
$ws.Range("L97").Value = 29000
$ws.Range("M97").Value = -3106.875
$ws.Range("N97").Value = -30982
$ws.Range("H101").Value = 24500
$ws.Range("J101").Value = 24500
$ws.Range("L101").Value = 24500
$ws.Range("N101").Value = -30990
$ws.Range("H103").Value = 32983.332
$ws.Range("J103").Value = 32983.332
$ws.Range("L103").Value = 32983.332
$ws.Range("N103").Value = -35327.332
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -28180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2368.6086
$ws.Range("I31").Value = 871.9459000000001
$ws.Range("K31").Value = 871.9459000000001
$ws.Range("M31").Value = -576.9459000000001
$ws.Range("H34").Value = 2368.6086
$ws.Range("I34").Value = 871.9459000000001
$ws.Range("K34").Value = 871.9459000000001
$ws.Range("M34").Value = -669.9459000000001
$ws.Range("H96").Value = 12979.6
$ws.Range("J96").Value = 12979.6
$ws.Range("L96").Value = 12979.6
$ws.Range("N96").Value = -18471.6
$ws.Range("H134").Value = 2166.375
$ws.Range("I134").Value = 1207.2307
$ws.Range("K134").Value = 3621.6921
$ws.Range("M134").Value = -1086.6921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 51.333332
$ws.Range("I38").Value = 26.428572
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 79.28571599999999
$ws.Range("L38").Value = 1200
$ws.Range("M38").Value = 267.714284
$ws.Range("N38").Value = -1894
$ws.Range("H86").Value = 614.9
$ws.Range("I86").Value = 631.1875
$ws.Range("J86").Value = 549.75
$ws.Range("K86").Value = 1893.5625
$ws.Range("L86").Value = 1649.25
$ws.Range("M86").Value = -707.5625
$ws.Range("N86").Value = -4021.25
$ws.Range("H89").Value = 614.9
$ws.Range("I89").Value = 631.1875
$ws.Range("J89").Value = 549.75
$ws.Range("K89").Value = 5680.6875
$ws.Range("L89").Value = 4947.75
$ws.Range("M89").Value = 247.3125
$ws.Range("N89").Value = -16803.75
$ws.Range("H100").Value = 6870.923
$ws.Range("J100").Value = 7108.0835
$ws.Range("L100").Value = 21324.2505
$ws.Range("N100").Value = -22946.2505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 26874.875
$ws.Range("J98").Value = 26874.875
$ws.Range("L98").Value = 26874.875
$ws.Range("N98").Value = -32864.875
$ws.Range("H99").Value = 7610.143
$ws.Range("I99").Value = 3654.2
$ws.Range("J99").Value = 17500
$ws.Range("K99").Value = 3654.2
$ws.Range("L99").Value = 17500
$ws.Range("M99").Value = -1408.2
$ws.Range("N99").Value = -21992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1554.4445
$ws.Range("I46").Value = 1252.1052
$ws.Range("J46").Value = 2272.5
$ws.Range("K46").Value = 1252.1052
$ws.Range("L46").Value = 2272.5
$ws.Range("M46").Value = -1064.1052
$ws.Range("N46").Value = -2648.5
$ws.Range("H100").Value = 47622884
$ws.Range("I100").Value = 4337.5
$ws.Range("J100").Value = 200002240
$ws.Range("K100").Value = 4337.5
$ws.Range("L100").Value = 200002240
$ws.Range("M100").Value = -3796.5
$ws.Range("N100").Value = -200003322
$ws.Range("H136").Value = 4316.095
$ws.Range("I136").Value = 2601.75
$ws.Range("J136").Value = 6601.8887
$ws.Range("K136").Value = 7805.25
$ws.Range("L136").Value = 19805.6661
$ws.Range("M136").Value = -5255.25
$ws.Range("N136").Value = -24905.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4125.263
$ws.Range("I62").Value = 3420
$ws.Range("J62").Value = 4313.3335
$ws.Range("K62").Value = 3420
$ws.Range("L62").Value = 4313.3335
$ws.Range("M62").Value = -2796
$ws.Range("N62").Value = -5561.3335
$ws.Range("H65").Value = 4125.263
$ws.Range("I65").Value = 3420
$ws.Range("J65").Value = 4313.3335
$ws.Range("K65").Value = 17100
$ws.Range("L65").Value = 21566.6675
$ws.Range("M65").Value = -13980
$ws.Range("N65").Value = -27806.6675
$ws.Range("H113").Value = 352.69565
$ws.Range("I113").Value = 303.94116
$ws.Range("J113").Value = 490.83334
$ws.Range("K113").Value = 911.82348
$ws.Range("L113").Value = 1472.50002
$ws.Range("M113").Value = 1258.17652
$ws.Range("N113").Value = -5812.500019999999
$ws.Range("H132").Value = 31252150
$ws.Range("I132").Value = 47620620
$ws.Range("J132").Value = 3253.3635
$ws.Range("K132").Value = 142861860
$ws.Range("L132").Value = 9760.0905
$ws.Range("M132").Value = -142859330
$ws.Range("N132").Value = -14820.0905
